$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1844155844155844
$ws.Range("C2").Value = 0.5701298701298702
$ws.Range("J2").Value = 0.02467532467532468
$ws.Range("P2").Value = 0.1454545454545454
$ws.Range("S2").Value = 0.07532467532467532

# Row 3
$ws.Range("B3").Value = 0.006437768240343348
$ws.Range("C3").Value = 0.03862660944206009
$ws.Range("J3").Value = 0.03004291845493562
$ws.Range("P3").Value = 0.7489270386266095
$ws.Range("S3").Value = 0.1759656652360515

# Row 4
$ws.Range("J4").Value = 0.03937007874015748
$ws.Range("P4").Value = 0.6692913385826772
$ws.Range("S4").Value = 0.2913385826771653

# Row 6
$ws.Range("B6").Value = 0.04748201438848921
$ws.Range("D6").Value = 0.01151079136690648
$ws.Range("F6").Value = 0.06474820143884892
$ws.Range("J6").Value = 0.2258992805755396
$ws.Range("O6").Value = 0.03165467625899281
$ws.Range("Q6").Value = 0.1640287769784173
$ws.Range("R6").Value = 0.08345323741007195
$ws.Range("S6").Value = 0.3712230215827338

# Row 7
$ws.Range("B7").Value = 0.08892921960072596
$ws.Range("D7").Value = 0.01633393829401089
$ws.Range("E7").Value = 0.003629764065335753
$ws.Range("F7").Value = 0.06533575317604355
$ws.Range("J7").Value = 0.132486388384755
$ws.Range("O7").Value = 0.03266787658802178
$ws.Range("Q7").Value = 0.1742286751361161
$ws.Range("R7").Value = 0.08348457350272233
$ws.Range("S7").Value = 0.4029038112522686

# Row 8
$ws.Range("B8").Value = 0.08096280087527352
$ws.Range("D8").Value = 0.0175054704595186
$ws.Range("F8").Value = 0.06564551422319474
$ws.Range("J8").Value = 0.1312910284463895
$ws.Range("O8").Value = 0.02334062727935813
$ws.Range("Q8").Value = 0.1684901531728665
$ws.Range("R8").Value = 0.1021152443471918
$ws.Range("S8").Value = 0.4106491611962071

# Row 9
$ws.Range("B9").Value = 0.08414239482200647
$ws.Range("D9").Value = 0.01456310679611651
$ws.Range("E9").Value = 0.001618122977346278
$ws.Range("F9").Value = 0.06472491909385113
$ws.Range("J9").Value = 0.1148867313915858
$ws.Range("O9").Value = 0.0226537216828479
$ws.Range("Q9").Value = 0.1666666666666667
$ws.Range("R9").Value = 0.1019417475728155
$ws.Range("S9").Value = 0.4288025889967638

# Row 10
$ws.Range("B10").Value = 0.09178255372945639
$ws.Range("D10").Value = 0.01946902654867257
$ws.Range("E10").Value = 0.001264222503160556
$ws.Range("F10").Value = 0.07281921618204804
$ws.Range("J10").Value = 0.1319848293299621
$ws.Range("O10").Value = 0.01769911504424779
$ws.Range("Q10").Value = 0.229582806573957
$ws.Range("R10").Value = 0.08874841972187104
$ws.Range("S10").Value = 0.3466498103666245

# Row 11
$ws.Range("G11").Value = 0.1425287356321839
$ws.Range("J11").Value = 0.0896551724137931
$ws.Range("K11").Value = 0.1896551724137931
$ws.Range("L11").Value = 0.5689655172413793
$ws.Range("S11").Value = 0.009195402298850575

# Row 12
$ws.Range("G12").Value = 0.7408829174664108
$ws.Range("J12").Value = 0.1880998080614203
$ws.Range("K12").Value = 0.007677543186180422
$ws.Range("L12").Value = 0.04606525911708254
$ws.Range("S12").Value = 0.01727447216890595

# Row 13
$ws.Range("G13").Value = 0.6063829787234043
$ws.Range("J13").Value = 0.3404255319148936
$ws.Range("S13").Value = 0.05319148936170213

# Row 15
$ws.Range("F15").Value = 0.01566951566951567
$ws.Range("H15").Value = 0.1566951566951567
$ws.Range("I15").Value = 0.0584045584045584
$ws.Range("J15").Value = 0.3774928774928775
$ws.Range("K15").Value = 0.05413105413105413
$ws.Range("M15").Value = 0.008547008547008548
$ws.Range("O15").Value = 0.0698005698005698
$ws.Range("S15").Value = 0.2592592592592592

# Row 16
$ws.Range("F16").Value = 0.01912045889101338
$ws.Range("H16").Value = 0.1720841300191205
$ws.Range("I16").Value = 0.0994263862332696
$ws.Range("J16").Value = 0.4225621414913958
$ws.Range("K16").Value = 0.09560229445506692
$ws.Range("M16").Value = 0.01720841300191205
$ws.Range("N16").Value = 0.001912045889101338
$ws.Range("O16").Value = 0.0497131931166348
$ws.Range("S16").Value = 0.1223709369024857

# Row 17
$ws.Range("F17").Value = 0.01246537396121884
$ws.Range("H17").Value = 0.1717451523545706
$ws.Range("I17").Value = 0.100415512465374
$ws.Range("J17").Value = 0.4411357340720222
$ws.Range("K17").Value = 0.09626038781163435
$ws.Range("M17").Value = 0.01246537396121884
$ws.Range("N17").Value = 0.002077562326869806
$ws.Range("O17").Value = 0.06163434903047092
$ws.Range("S17").Value = 0.1018005540166205

# Row 18
$ws.Range("F18").Value = 0.0245398773006135
$ws.Range("H18").Value = 0.1886503067484663
$ws.Range("I18").Value = 0.07822085889570553
$ws.Range("J18").Value = 0.4217791411042945
$ws.Range("K18").Value = 0.1012269938650307
$ws.Range("M18").Value = 0.01226993865030675
$ws.Range("N18").Value = 0.001533742331288344
$ws.Range("O18").Value = 0.07515337423312883
$ws.Range("S18").Value = 0.09662576687116564

# Row 19
$ws.Range("F19").Value = 0.01474201474201474
$ws.Range("H19").Value = 0.2148512148512149
$ws.Range("I19").Value = 0.09009009009009009
$ws.Range("J19").Value = 0.3808353808353808
$ws.Range("K19").Value = 0.1094731094731095
$ws.Range("M19").Value = 0.01528801528801529
$ws.Range("N19").Value = 0.000546000546000546
$ws.Range("O19").Value = 0.07316407316407317
$ws.Range("S19").Value = 0.101010101010101
